$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.8.1 -> 1.8.2
$ws.Range("B3").Value = "1.8.2"

# Experimental value cleared (was "true")
$ws.Range("B7").Value = ""

# Date: 2023-10-31 -> 2025-11-18
# (written via a scratch cell + paste-values so the literal text is stored
# as a shared string without Excel's autodetect turning it into a date
# serial number / introducing a new number format)
$ws.Range("Z100").Formula = "=""2025-11-18"""
$ws.Range("Z100").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("Z100").Clear()
